$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A7").Value = "I1_S7"
$ws.Range("B7").Value = "rispetto alla precedente penalità 1"

$ws.Range("A8").Value = "I1_S8"
$ws.Range("B8").Value = "rispetto alla precedente penalità 5"

$ws.Range("B9").Select()
